# Auto-generated: apply cached-value updates to Leve profit columns (H-N)
# across all 8 job sheets, per the scheduled market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 219.22223   # H5: 229.375 -> 219.22223
$ws.Cells.Item(5, 9).Value = 103.42857   # I5: 97.666664 -> 103.42857
$ws.Cells.Item(5, 11).Value = 103.42857   # K5: 97.666664 -> 103.42857
$ws.Cells.Item(5, 13).Value = 11.57143000000001   # M5: 17.333336 -> 11.57143000000001
$ws.Cells.Item(43, 8).Value = 4447.5   # H43: 4966.3335 -> 4447.5
$ws.Cells.Item(43, 9).Value = 4965.3335   # I43: 4966.3335 -> 4965.3335
$ws.Cells.Item(43, 10).Value = 2894   # J43: 0 -> 2894
$ws.Cells.Item(43, 11).Value = 4965.3335   # K43: 4966.3335 -> 4965.3335
$ws.Cells.Item(43, 12).Value = 2894   # L43: 0 -> 2894
$ws.Cells.Item(43, 13).Value = -4896.3335   # M43: -4897.3335 -> -4896.3335
$ws.Cells.Item(43, 14).Value = -3032   # N43: None -> -3032
$ws.Cells.Item(51, 8).Value = 7109.3335   # H51: 7373.375 -> 7109.3335
$ws.Cells.Item(51, 9).Value = 6498   # I51: 6712.4287 -> 6498
$ws.Cells.Item(51, 11).Value = 6498   # K51: 6712.4287 -> 6498
$ws.Cells.Item(51, 13).Value = -6014   # M51: -6228.4287 -> -6014
$ws.Cells.Item(55, 8).Value = 278.14285   # H55: 264.875 -> 278.14285
$ws.Cells.Item(55, 9).Value = 178   # I55: 176 -> 178
$ws.Cells.Item(55, 11).Value = 178   # K55: 176 -> 178
$ws.Cells.Item(55, 13).Value = 36   # M55: 38 -> 36
$ws.Cells.Item(80, 8).Value = 3060.818   # H80: 4125 -> 3060.818
$ws.Cells.Item(80, 9).Value = 2639.8   # I80: 3250 -> 2639.8
$ws.Cells.Item(80, 10).Value = 3411.6667   # J80: 5000 -> 3411.6667
$ws.Cells.Item(80, 11).Value = 7919.400000000001   # K80: 9750 -> 7919.400000000001
$ws.Cells.Item(80, 12).Value = 10235.0001   # L80: 15000 -> 10235.0001
$ws.Cells.Item(80, 13).Value = -6921.400000000001   # M80: -8752 -> -6921.400000000001
$ws.Cells.Item(80, 14).Value = -12231.0001   # N80: -16996 -> -12231.0001
$ws.Cells.Item(83, 8).Value = 3060.818   # H83: 4125 -> 3060.818
$ws.Cells.Item(83, 9).Value = 2639.8   # I83: 3250 -> 2639.8
$ws.Cells.Item(83, 10).Value = 3411.6667   # J83: 5000 -> 3411.6667
$ws.Cells.Item(83, 11).Value = 23758.2   # K83: 29250 -> 23758.2
$ws.Cells.Item(83, 12).Value = 30705.0003   # L83: 45000 -> 30705.0003
$ws.Cells.Item(83, 13).Value = -18766.2   # M83: -24258 -> -18766.2
$ws.Cells.Item(83, 14).Value = -40689.0003   # N83: -54984 -> -40689.0003
$ws.Cells.Item(86, 8).Value = 6569.778   # H86: 5556.091 -> 6569.778
$ws.Cells.Item(86, 9).Value = 5450   # I86: 4558.8 -> 5450
$ws.Cells.Item(86, 10).Value = 7465.6   # J86: 6387.1665 -> 7465.6
$ws.Cells.Item(86, 11).Value = 5450   # K86: 4558.8 -> 5450
$ws.Cells.Item(86, 12).Value = 7465.6   # L86: 6387.1665 -> 7465.6
$ws.Cells.Item(86, 13).Value = -4327   # M86: -3435.8 -> -4327
$ws.Cells.Item(86, 14).Value = -9711.6   # N86: -8633.166499999999 -> -9711.6
$ws.Cells.Item(89, 8).Value = 6569.778   # H89: 5556.091 -> 6569.778
$ws.Cells.Item(89, 9).Value = 5450   # I89: 4558.8 -> 5450
$ws.Cells.Item(89, 10).Value = 7465.6   # J89: 6387.1665 -> 7465.6
$ws.Cells.Item(89, 11).Value = 27250   # K89: 22794 -> 27250
$ws.Cells.Item(89, 12).Value = 37328   # L89: 31935.8325 -> 37328
$ws.Cells.Item(89, 13).Value = -21634   # M89: -17178 -> -21634
$ws.Cells.Item(89, 14).Value = -48560   # N89: -43167.8325 -> -48560
$ws.Cells.Item(135, 8).Value = 1244.8572   # H135: 1596.8572 -> 1244.8572
$ws.Cells.Item(135, 9).Value = 235.2   # I135: 254.55556 -> 235.2
$ws.Cells.Item(135, 10).Value = 3769   # J135: 4013 -> 3769
$ws.Cells.Item(135, 11).Value = 2116.8   # K135: 2291.00004 -> 2116.8
$ws.Cells.Item(135, 12).Value = 33921   # L135: 36117 -> 33921
$ws.Cells.Item(135, 13).Value = 418.2000000000003   # M135: 243.9999600000001 -> 418.2000000000003
$ws.Cells.Item(135, 14).Value = -38991   # N135: -41187 -> -38991
$ws.Cells.Item(138, 8).Value = 2542.75   # H138: 2651.8235 -> 2542.75
$ws.Cells.Item(138, 9).Value = 1910.875   # I138: 2187.111 -> 1910.875
$ws.Cells.Item(138, 11).Value = 5732.625   # K138: 6561.333 -> 5732.625
$ws.Cells.Item(138, 13).Value = -592.625   # M138: -1421.333 -> -592.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 3347   # H5: 2521.5 -> 3347
$ws.Cells.Item(5, 10).Value = 10000   # J5: 5022.5 -> 10000
$ws.Cells.Item(5, 12).Value = 10000   # L5: 5022.5 -> 10000
$ws.Cells.Item(5, 14).Value = -10224   # N5: -5246.5 -> -10224
$ws.Cells.Item(10, 8).Value = 7500   # H10: 10005 -> 7500
$ws.Cells.Item(10, 10).Value = 7500   # J10: 10005 -> 7500
$ws.Cells.Item(10, 12).Value = 7500   # L10: 10005 -> 7500
$ws.Cells.Item(10, 14).Value = -7840   # N10: -10345 -> -7840
$ws.Cells.Item(24, 8).Value = 67049.75   # H24: 56592.332 -> 67049.75
$ws.Cells.Item(24, 10).Value = 67049.75   # J24: 56592.332 -> 67049.75
$ws.Cells.Item(24, 12).Value = 67049.75   # L24: 56592.332 -> 67049.75
$ws.Cells.Item(24, 14).Value = -67797.75   # N24: -57340.332 -> -67797.75
$ws.Cells.Item(96, 8).Value = 20000   # H96: 15000 -> 20000
$ws.Cells.Item(96, 10).Value = 20000   # J96: 15000 -> 20000
$ws.Cells.Item(96, 12).Value = 20000   # L96: 15000 -> 20000
$ws.Cells.Item(96, 14).Value = -25492   # N96: -20492 -> -25492
$ws.Cells.Item(100, 8).Value = 67049.75   # H100: 56592.332 -> 67049.75
$ws.Cells.Item(100, 10).Value = 67049.75   # J100: 56592.332 -> 67049.75
$ws.Cells.Item(100, 12).Value = 67049.75   # L100: 56592.332 -> 67049.75
$ws.Cells.Item(100, 14).Value = -69213.75   # N100: -58756.332 -> -69213.75
$ws.Cells.Item(139, 8).Value = 87500   # H139: 55000 -> 87500
$ws.Cells.Item(139, 10).Value = 87500   # J139: 55000 -> 87500
$ws.Cells.Item(139, 12).Value = 87500   # L139: 55000 -> 87500
$ws.Cells.Item(139, 14).Value = -97780   # N139: -65280 -> -97780
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 3347   # H4: 2521.5 -> 3347
$ws.Cells.Item(4, 10).Value = 10000   # J4: 5022.5 -> 10000
$ws.Cells.Item(4, 12).Value = 10000   # L4: 5022.5 -> 10000
$ws.Cells.Item(4, 14).Value = -10230   # N4: -5252.5 -> -10230
$ws.Cells.Item(22, 8).Value = 1333.3334   # H22: 714.2857 -> 1333.3334
$ws.Cells.Item(22, 9).Value = 1333.3334   # I22: 800 -> 1333.3334
$ws.Cells.Item(22, 10).Value = 0   # J22: 200 -> 0
$ws.Cells.Item(22, 11).Value = 1333.3334   # K22: 800 -> 1333.3334
$ws.Cells.Item(22, 12).Value = 0   # L22: 200 -> 0
$ws.Cells.Item(22, 13).ClearContents()   # M22: clear (was -627)
$ws.Cells.Item(22, 14).Value = -1160.3334   # N22: -546 -> -1160.3334
$ws.Cells.Item(29, 8).Value = 515   # H29: 507.5 -> 515
$ws.Cells.Item(29, 9).Value = 515   # I29: 507.5 -> 515
$ws.Cells.Item(29, 11).Value = 515   # K29: 507.5 -> 515
$ws.Cells.Item(29, 13).Value = -226   # M29: -218.5 -> -226
$ws.Cells.Item(107, 8).Value = 939.2857   # H107: 1065.8334 -> 939.2857
$ws.Cells.Item(107, 9).Value = 735   # I107: 873.75 -> 735
$ws.Cells.Item(107, 11).Value = 735   # K107: 873.75 -> 735
$ws.Cells.Item(107, 13).Value = 1185   # M107: 1046.25 -> 1185
$ws.Cells.Item(134, 8).Value = 2826.077   # H134: 3374.5 -> 2826.077
$ws.Cells.Item(134, 9).Value = 2826.077   # I134: 3374.5 -> 2826.077
$ws.Cells.Item(134, 11).Value = 8478.231   # K134: 10123.5 -> 8478.231
$ws.Cells.Item(134, 13).Value = -5943.231   # M134: -7588.5 -> -5943.231
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2499.1667   # H31: 2836.75 -> 2499.1667
$ws.Cells.Item(31, 9).Value = 1882   # I31: 1998 -> 1882
$ws.Cells.Item(31, 11).Value = 1882   # K31: 1998 -> 1882
$ws.Cells.Item(31, 13).Value = -1587   # M31: -1703 -> -1587
$ws.Cells.Item(34, 8).Value = 2499.1667   # H34: 2836.75 -> 2499.1667
$ws.Cells.Item(34, 9).Value = 1882   # I34: 1998 -> 1882
$ws.Cells.Item(34, 11).Value = 1882   # K34: 1998 -> 1882
$ws.Cells.Item(34, 13).Value = -1680   # M34: -1796 -> -1680
$ws.Cells.Item(106, 8).Value = 24214.143   # H106: 92428.42999999999 -> 24214.143
$ws.Cells.Item(106, 10).Value = 24214.143   # J106: 92428.42999999999 -> 24214.143
$ws.Cells.Item(106, 12).Value = 24214.143   # L106: 92428.42999999999 -> 24214.143
$ws.Cells.Item(106, 14).Value = -26738.143   # N106: -94952.42999999999 -> -26738.143
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 150.66667   # H38: 173.3 -> 150.66667
$ws.Cells.Item(38, 9).Value = 61.42857   # I38: 69.25 -> 61.42857
$ws.Cells.Item(38, 10).Value = 275.6   # J38: 242.66667 -> 275.6
$ws.Cells.Item(38, 11).Value = 184.28571   # K38: 207.75 -> 184.28571
$ws.Cells.Item(38, 12).Value = 826.8000000000001   # L38: 728.00001 -> 826.8000000000001
$ws.Cells.Item(38, 13).Value = 162.71429   # M38: 139.25 -> 162.71429
$ws.Cells.Item(38, 14).Value = -1520.8   # N38: -1422.00001 -> -1520.8
$ws.Cells.Item(105, 8).Value = 0   # H105: 7000 -> 0
$ws.Cells.Item(105, 9).Value = 0   # I105: 7000 -> 0
$ws.Cells.Item(105, 11).Value = 0   # K105: 21000 -> 0
$ws.Cells.Item(105, 13).ClearContents()   # M105: clear (was -18379)
$ws.Cells.Item(126, 8).Value = 2500   # H126: 0 -> 2500
$ws.Cells.Item(126, 9).Value = 2500   # I126: 0 -> 2500
$ws.Cells.Item(126, 11).Value = 7500   # K126: 0 -> 7500
$ws.Cells.Item(126, 13).Value = -2560   # M126: None -> -2560
$ws.Cells.Item(130, 8).Value = 0   # H130: 17250 -> 0
$ws.Cells.Item(130, 9).Value = 0   # I130: 14500 -> 0
$ws.Cells.Item(130, 10).Value = 0   # J130: 20000 -> 0
$ws.Cells.Item(130, 11).Value = 0   # K130: 43500 -> 0
$ws.Cells.Item(130, 12).ClearContents()   # L130: clear (was 60000)
$ws.Cells.Item(130, 13).ClearContents()   # M130: clear (was -38480)
$ws.Cells.Item(130, 14).Value = 0   # N130: -70040 -> 0
$ws.Cells.Item(131, 8).Value = 4134   # H131: 1579.7693 -> 4134
$ws.Cells.Item(131, 10).Value = 4134   # J131: 1579.7693 -> 4134
$ws.Cells.Item(131, 12).Value = 12402   # L131: 4739.3079 -> 12402
$ws.Cells.Item(131, 14).Value = -22482   # N131: -14819.3079 -> -22482
$ws.Cells.Item(134, 8).Value = 0   # H134: 2874.75 -> 0
$ws.Cells.Item(134, 9).Value = 0   # I134: 2874.75 -> 0
$ws.Cells.Item(134, 11).Value = 0   # K134: 8624.25 -> 0
$ws.Cells.Item(134, 13).ClearContents()   # M134: clear (was -3554.25)
$ws.Cells.Item(136, 8).Value = 0   # H136: 3376.6667 -> 0
$ws.Cells.Item(136, 9).Value = 0   # I136: 3376.6667 -> 0
$ws.Cells.Item(136, 11).Value = 0   # K136: 10130.0001 -> 0
$ws.Cells.Item(136, 13).ClearContents()   # M136: clear (was -5030.000100000001)
$ws.Cells.Item(137, 8).Value = 2521.889   # H137: 2339.5454 -> 2521.889
$ws.Cells.Item(137, 9).Value = 1357.25   # I137: 1492.5 -> 1357.25
$ws.Cells.Item(137, 10).Value = 3453.6   # J137: 2823.5715 -> 3453.6
$ws.Cells.Item(137, 11).Value = 4071.75   # K137: 4477.5 -> 4071.75
$ws.Cells.Item(137, 12).Value = 10360.8   # L137: 8470.7145 -> 10360.8
$ws.Cells.Item(137, 13).Value = 1028.25   # M137: 622.5 -> 1028.25
$ws.Cells.Item(137, 14).Value = -20560.8   # N137: -18670.7145 -> -20560.8
$ws.Cells.Item(138, 8).Value = 0   # H138: 3196.6 -> 0
$ws.Cells.Item(138, 9).Value = 0   # I138: 3245.75 -> 0
$ws.Cells.Item(138, 10).Value = 0   # J138: 3000 -> 0
$ws.Cells.Item(138, 11).Value = 0   # K138: 9737.25 -> 0
$ws.Cells.Item(138, 12).ClearContents()   # L138: clear (was 9000)
$ws.Cells.Item(138, 13).ClearContents()   # M138: clear (was -4597.25)
$ws.Cells.Item(138, 14).Value = 0   # N138: -19280 -> 0
$ws.Cells.Item(139, 8).Value = 10000   # H139: 2539.5 -> 10000
$ws.Cells.Item(139, 9).Value = 10000   # I139: 3398.75 -> 10000
$ws.Cells.Item(139, 10).Value = 0   # J139: 1966.6666 -> 0
$ws.Cells.Item(139, 11).Value = 30000   # K139: 10196.25 -> 30000
$ws.Cells.Item(139, 12).Value = 0   # L139: 5899.9998 -> 0
$ws.Cells.Item(139, 13).ClearContents()   # M139: clear (was -5056.25)
$ws.Cells.Item(139, 14).Value = -24860   # N139: -16179.9998 -> -24860
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 25000   # H15: 24999.666 -> 25000
$ws.Cells.Item(15, 10).Value = 25000   # J15: 24999.666 -> 25000
$ws.Cells.Item(15, 12).Value = 25000   # L15: 24999.666 -> 25000
$ws.Cells.Item(15, 14).Value = -25576   # N15: -25575.666 -> -25576
$ws.Cells.Item(80, 8).Value = 4734   # H80: 2908.7144 -> 4734
$ws.Cells.Item(80, 9).Value = 962   # I80: 712.25 -> 962
$ws.Cells.Item(80, 10).Value = 8506   # J80: 5837.3335 -> 8506
$ws.Cells.Item(80, 11).Value = 962   # K80: 712.25 -> 962
$ws.Cells.Item(80, 12).Value = 8506   # L80: 5837.3335 -> 8506
$ws.Cells.Item(80, 13).Value = 36   # M80: 285.75 -> 36
$ws.Cells.Item(80, 14).Value = -10502   # N80: -7833.3335 -> -10502
$ws.Cells.Item(81, 8).Value = 25000   # H81: 24999.666 -> 25000
$ws.Cells.Item(81, 10).Value = 25000   # J81: 24999.666 -> 25000
$ws.Cells.Item(81, 12).Value = 25000   # L81: 24999.666 -> 25000
$ws.Cells.Item(81, 14).Value = -26996   # N81: -26995.666 -> -26996
$ws.Cells.Item(83, 8).Value = 4734   # H83: 2908.7144 -> 4734
$ws.Cells.Item(83, 9).Value = 962   # I83: 712.25 -> 962
$ws.Cells.Item(83, 10).Value = 8506   # J83: 5837.3335 -> 8506
$ws.Cells.Item(83, 11).Value = 4810   # K83: 3561.25 -> 4810
$ws.Cells.Item(83, 12).Value = 42530   # L83: 29186.6675 -> 42530
$ws.Cells.Item(83, 13).Value = 182   # M83: 1430.75 -> 182
$ws.Cells.Item(83, 14).Value = -52514   # N83: -39170.6675 -> -52514
$ws.Cells.Item(84, 8).Value = 25000   # H84: 24999.666 -> 25000
$ws.Cells.Item(84, 10).Value = 25000   # J84: 24999.666 -> 25000
$ws.Cells.Item(84, 12).Value = 75000   # L84: 74998.99800000001 -> 75000
$ws.Cells.Item(84, 14).Value = -84984   # N84: -84982.99800000001 -> -84984
$ws.Cells.Item(98, 8).Value = 15299.8   # H98: 13126.25 -> 15299.8
$ws.Cells.Item(98, 10).Value = 15299.8   # J98: 13126.25 -> 15299.8
$ws.Cells.Item(98, 12).Value = 15299.8   # L98: 13126.25 -> 15299.8
$ws.Cells.Item(98, 14).Value = -21289.8   # N98: -19116.25 -> -21289.8
$ws.Cells.Item(122, 8).Value = 1678.5714   # H122: 1618.75 -> 1678.5714
$ws.Cells.Item(122, 9).Value = 1678.5714   # I122: 1618.75 -> 1678.5714
$ws.Cells.Item(122, 11).Value = 5035.7142   # K122: 4856.25 -> 5035.7142
$ws.Cells.Item(122, 13).Value = -2585.7142   # M122: -2406.25 -> -2585.7142
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3233   # H22: 3526.3 -> 3233
$ws.Cells.Item(22, 9).Value = 400   # I22: 0 -> 400
$ws.Cells.Item(22, 10).Value = 3862.5557   # J22: 3526.3 -> 3862.5557
$ws.Cells.Item(22, 11).Value = 400   # K22: 0 -> 400
$ws.Cells.Item(22, 12).Value = 3862.5557   # L22: 3526.3 -> 3862.5557
$ws.Cells.Item(22, 13).Value = -105   # M22: None -> -105
$ws.Cells.Item(22, 14).Value = -4452.5557   # N22: -4116.3 -> -4452.5557
$ws.Cells.Item(27, 8).Value = 3233   # H27: 3526.3 -> 3233
$ws.Cells.Item(27, 9).Value = 400   # I27: 0 -> 400
$ws.Cells.Item(27, 10).Value = 3862.5557   # J27: 3526.3 -> 3862.5557
$ws.Cells.Item(27, 11).Value = 400   # K27: 0 -> 400
$ws.Cells.Item(27, 12).Value = 3862.5557   # L27: 3526.3 -> 3862.5557
$ws.Cells.Item(27, 13).Value = -293   # M27: None -> -293
$ws.Cells.Item(27, 14).Value = -4076.5557   # N27: -3740.3 -> -4076.5557
$ws.Cells.Item(46, 8).Value = 4846.077   # H46: 4914.2144 -> 4846.077
$ws.Cells.Item(46, 10).Value = 5999.857   # J46: 5974.875 -> 5999.857
$ws.Cells.Item(46, 12).Value = 5999.857   # L46: 5974.875 -> 5999.857
$ws.Cells.Item(46, 14).Value = -6375.857   # N46: -6350.875 -> -6375.857
$ws.Cells.Item(99, 8).Value = 17078   # H99: 17717.5 -> 17078
$ws.Cells.Item(99, 9).Value = 17078   # I99: 17717.5 -> 17078
$ws.Cells.Item(99, 11).Value = 17078   # K99: 17717.5 -> 17078
$ws.Cells.Item(99, 13).Value = -14083   # M99: -14722.5 -> -14083
$ws.Cells.Item(101, 8).Value = 19430.5   # H101: 19575.666 -> 19430.5
$ws.Cells.Item(101, 10).Value = 19430.5   # J101: 19575.666 -> 19430.5
$ws.Cells.Item(101, 12).Value = 19430.5   # L101: 19575.666 -> 19430.5
$ws.Cells.Item(101, 14).Value = -25920.5   # N101: -26065.666 -> -25920.5
$ws.Cells.Item(130, 8).Value = 57149.5   # H130: 0 -> 57149.5
$ws.Cells.Item(130, 10).Value = 57149.5   # J130: 0 -> 57149.5
$ws.Cells.Item(130, 12).Value = 57149.5   # L130: 0 -> 57149.5
$ws.Cells.Item(130, 14).Value = -67189.5   # N130: None -> -67189.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(93, 8).Value = 24000   # H93: 17500 -> 24000
$ws.Cells.Item(93, 9).Value = 0   # I93: 11000 -> 0
$ws.Cells.Item(93, 11).Value = 0   # K93: 11000 -> 0
$ws.Cells.Item(93, 13).ClearContents()   # M93: clear (was -8504)
$ws.Cells.Item(119, 8).Value = 40999.668   # H119: 57666.332 -> 40999.668
$ws.Cells.Item(119, 10).Value = 40999.668   # J119: 57666.332 -> 40999.668
$ws.Cells.Item(119, 12).Value = 40999.668   # L119: 57666.332 -> 40999.668
$ws.Cells.Item(119, 14).Value = -50675.668   # N119: -67342.33199999999 -> -50675.668
$ws.Cells.Item(136, 8).Value = 1064   # H136: 1101.5555 -> 1064
$ws.Cells.Item(136, 9).Value = 1081.4445   # I136: 1122.2354 -> 1081.4445
$ws.Cells.Item(136, 11).Value = 3244.3335   # K136: 3366.7062 -> 3244.3335
$ws.Cells.Item(136, 13).Value = -694.3335000000002   # M136: -816.7062000000001 -> -694.3335000000002
